$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 224.72728
$ws.Range("I4").Value = 147.2
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 147.2
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -33.19999999999999
$ws.Range("N4").Value = -1228
# Row 12
$ws.Range("H12").Value = 1324.75
$ws.Range("I12").Value = 1150.5
$ws.Range("K12").Value = 1150.5
$ws.Range("M12").Value = -980.5
# Row 33
$ws.Range("H33").Value = 88.4
$ws.Range("J33").Value = 49.5
$ws.Range("L33").Value = 49.5
$ws.Range("N33").Value = -507.5
# Row 53
$ws.Range("H53").Value = 708.8182
$ws.Range("I53").Value = 868.1667
$ws.Range("K53").Value = 868.1667
$ws.Range("M53").Value = -231.1667
# Row 112
$ws.Range("H112").Value = 4156.852
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 4274.423
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 12823.269
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -15039.269
# Row 132
$ws.Range("H132").Value = 712.5217
$ws.Range("I132").Value = 738.1429
$ws.Range("J132").Value = 443.5
$ws.Range("K132").Value = 2214.4287
$ws.Range("L132").Value = 1330.5
$ws.Range("M132").Value = 315.5712999999996
$ws.Range("N132").Value = -6390.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3037513
$ws.Range("I32").Value = 3284112.5
$ws.Range("K32").Value = 3284112.5
$ws.Range("M32").Value = -3283825.5
# Row 45
$ws.Range("H45").Value = 6215.4375
$ws.Range("I45").Value = 3500
$ws.Range("K45").Value = 3500
$ws.Range("M45").Value = -3123
# Row 57
$ws.Range("H57").Value = 5499.4443
$ws.Range("I57").Value = 5499.4443
$ws.Range("K57").Value = 5499.4443
$ws.Range("M57").Value = -5015.4443
# Row 61
$ws.Range("H61").Value = 50007028
$ws.Range("I61").Value = 1713.5834
$ws.Range("J61").Value = 125015000
$ws.Range("K61").Value = 1713.5834
$ws.Range("L61").Value = 125015000
$ws.Range("M61").Value = -1501.5834
$ws.Range("N61").Value = -125015424
# Row 74
$ws.Range("H74").Value = 45118.418
$ws.Range("I74").Value = 73603.43
$ws.Range("K74").Value = 73603.43
$ws.Range("M74").Value = -72729.43
# Row 77
$ws.Range("H77").Value = 45118.418
$ws.Range("I77").Value = 73603.43
$ws.Range("K77").Value = 368017.15
$ws.Range("M77").Value = -363649.15
# Row 106
$ws.Range("H106").Value = 46894.332
$ws.Range("J106").Value = 46894.332
$ws.Range("L106").Value = 46894.332
$ws.Range("N106").Value = -49418.332
# Row 108
$ws.Range("H108").Value = 59376
$ws.Range("J108").Value = 59376
$ws.Range("L108").Value = 59376
$ws.Range("N108").Value = -67056
# Row 110
$ws.Range("H110").Value = 111112250
$ws.Range("I110").Value = 1694.5
$ws.Range("K110").Value = 1694.5
$ws.Range("M110").Value = 350.5
# Row 115
$ws.Range("H115").Value = 59376
$ws.Range("J115").Value = 59376
$ws.Range("L115").Value = 59376
$ws.Range("N115").Value = -62510
# Row 122
$ws.Range("H122").Value = 5667.6
$ws.Range("I122").Value = 3155.5
$ws.Range("K122").Value = 9466.5
$ws.Range("M122").Value = -7016.5
# Row 126
$ws.Range("H126").Value = 5397.5713
$ws.Range("I126").Value = 5397.5713
$ws.Range("K126").Value = 16192.7139
$ws.Range("M126").Value = -13722.7139
# Row 132
$ws.Range("H132").Value = 5811
$ws.Range("I132").Value = 1980.65
$ws.Range("J132").Value = 10317.294
$ws.Range("K132").Value = 5941.950000000001
$ws.Range("L132").Value = 30951.882
$ws.Range("M132").Value = -3411.950000000001
$ws.Range("N132").Value = -36011.882
# Row 136
$ws.Range("H136").Value = 50007028
$ws.Range("I136").Value = 1713.5834
$ws.Range("J136").Value = 125015000
$ws.Range("K136").Value = 5140.7502
$ws.Range("L136").Value = 375045000
$ws.Range("M136").Value = -2590.7502
$ws.Range("N136").Value = -375050100

$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 516
$ws.Range("I29").Value = 516
$ws.Range("K29").Value = 516
$ws.Range("M29").Value = -227
# Row 30
$ws.Range("H30").Value = 1140
$ws.Range("J30").Value = 1140
$ws.Range("L30").Value = 1140
$ws.Range("N30").Value = -1390
# Row 105
$ws.Range("H105").Value = 4458.3076
$ws.Range("I105").Value = 3658
$ws.Range("K105").Value = 3658
$ws.Range("M105").Value = -1911
# Row 128
$ws.Range("H128").Value = 4512.222
$ws.Range("I128").Value = 4512.222
$ws.Range("K128").Value = 13536.666
$ws.Range("M128").Value = -11046.666
# Row 134
$ws.Range("H134").Value = 15635531
$ws.Range("I134").Value = 125003000
$ws.Range("J134").Value = 11607.143
$ws.Range("K134").Value = 375009000
$ws.Range("L134").Value = 34821.429
$ws.Range("M134").Value = -375006465
$ws.Range("N134").Value = -39891.429

$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Range("H21").Value = 10960.333
$ws.Range("J21").Value = 10960.333
$ws.Range("L21").Value = 10960.333
$ws.Range("N21").Value = -11430.333
# Row 22
$ws.Range("H22").Value = 824.75
$ws.Range("I22").Value = 933.3333
$ws.Range("J22").Value = 499
$ws.Range("K22").Value = 933.3333
$ws.Range("L22").Value = 499
$ws.Range("M22").Value = -583.3333
$ws.Range("N22").Value = -1199
# Row 43
$ws.Range("H43").Value = 32868
$ws.Range("J43").Value = 32868
$ws.Range("L43").Value = 32868
$ws.Range("N43").Value = -33236
# Row 76
$ws.Range("H76").Value = 5299.1665
$ws.Range("I76").Value = 5299.1665
$ws.Range("K76").Value = 5299.1665
$ws.Range("M76").Value = -4984.1665
# Row 79
$ws.Range("H79").Value = 5299.1665
$ws.Range("I79").Value = 5299.1665
$ws.Range("K79").Value = 5299.1665
$ws.Range("M79").Value = -4207.1665
# Row 86
$ws.Range("H86").Value = 17366458
$ws.Range("I86").Value = 24044412
$ws.Range("K86").Value = 24044412
$ws.Range("M86").Value = -24043289
# Row 89
$ws.Range("H89").Value = 17366458
$ws.Range("I89").Value = 24044412
$ws.Range("K89").Value = 120222060
$ws.Range("M89").Value = -120216444
# Row 101
$ws.Range("H101").Value = 32868
$ws.Range("J101").Value = 32868
$ws.Range("L101").Value = 32868
$ws.Range("N101").Value = -39358
# Row 134
$ws.Range("H134").Value = 4533.186
$ws.Range("I134").Value = 1662.9667
$ws.Range("K134").Value = 4988.9001
$ws.Range("M134").Value = -2453.9001

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 23809852
$ws.Range("I33").Value = 55555640
$ws.Range("J33").Value = 513.5
$ws.Range("K33").Value = 333333840
$ws.Range("L33").Value = 3081
$ws.Range("M33").Value = -333333557
$ws.Range("N33").Value = -3647
# Row 131
$ws.Range("H131").Value = 37335.645
$ws.Range("I131").Value = 2333.3333
$ws.Range("J131").Value = 41535.92
$ws.Range("K131").Value = 6999.999899999999
$ws.Range("L131").Value = 124607.76
$ws.Range("M131").Value = -1959.999899999999
$ws.Range("N131").Value = -134687.76

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4531.394
$ws.Range("I132").Value = 1775.0476
$ws.Range("J132").Value = 9355
$ws.Range("K132").Value = 5325.142800000001
$ws.Range("L132").Value = 28065
$ws.Range("M132").Value = -2795.142800000001
$ws.Range("N132").Value = -33125
# Row 136
$ws.Range("H136").Value = 27754.6
$ws.Range("J136").Value = 28145.52
$ws.Range("L136").Value = 84436.56
$ws.Range("N136").Value = -89536.56

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 2837.8572
$ws.Range("I100").Value = 2894.3333
$ws.Range("J100").Value = 2499
$ws.Range("K100").Value = 2894.3333
$ws.Range("L100").Value = 2499
$ws.Range("M100").Value = -2353.3333
$ws.Range("N100").Value = -3581
# Row 107
$ws.Range("H107").Value = 4570.5713
$ws.Range("I107").Value = 4570.5713
$ws.Range("K107").Value = 4570.5713
$ws.Range("M107").Value = -2650.5713
# Row 132
$ws.Range("H132").Value = 9809725
$ws.Range("I132").Value = 17243964
$ws.Range("K132").Value = 51731892
$ws.Range("M132").Value = -51729362
# Row 136
$ws.Range("H136").Value = 11504.106
$ws.Range("I136").Value = 2647.8
$ws.Range("J136").Value = 21568.092
$ws.Range("K136").Value = 7943.400000000001
$ws.Range("L136").Value = 64704.276
$ws.Range("M136").Value = -5393.400000000001
$ws.Range("N136").Value = -69804.276

